# Updated cryptos list — applies per-cell price/volume(1h) refresh
# plus the RocketPoolETH / TrustWalletToken row swap (rows 45-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.542.91"
$ws.Range("E2").Value = "'  -3.01%  "
$ws.Range("D3").Value = "'1.660.81"
$ws.Range("E3").Value = "'  -3.68%  "
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("D5").Value = "'214.55"
$ws.Range("E5").Value = "'  -1.81%  "
$ws.Range("D6").Value = "'0.512"
$ws.Range("E6").Value = "'  -2.08%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D8").Value = "'24.38"
$ws.Range("E8").Value = "'  +1.94%  "
$ws.Range("E9").Value = "'  -1.58%  "
$ws.Range("D10").Value = "'0.0619"
$ws.Range("E10").Value = "'  -2.36%  "
$ws.Range("D11").Value = "'0.0878"
$ws.Range("E11").Value = "'  -1.65%  "
$ws.Range("D12").Value = "'1.896.02"
$ws.Range("E12").Value = "'  -3.73%  "
$ws.Range("D13").Value = "'1.661.54"
$ws.Range("E13").Value = "'  -3.68%  "
$ws.Range("D14").Value = "'4.13"
$ws.Range("E14").Value = "'  -2.29%  "
$ws.Range("D15").Value = "'0.566"
$ws.Range("E15").Value = "'  +0.66%  "
$ws.Range("D16").Value = "'65.90"
$ws.Range("E16").Value = "'  -2.45%  "
$ws.Range("D17").Value = "'27.542.25"
$ws.Range("D18").Value = "'239.95"
$ws.Range("E18").Value = "'  -2.73%  "
$ws.Range("D19").Value = "'0.0₃0729"
$ws.Range("E19").Value = "'  -2.71%  "
$ws.Range("E20").Value = "'  -2.49%  "
$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("D22").Value = "'4.45"
$ws.Range("E22").Value = "'  -3.18%  "
$ws.Range("D23").Value = "'9.40"
$ws.Range("E23").Value = "'  -2.36%  "
$ws.Range("E24").Value = "'  -1.38%  "
$ws.Range("D25").Value = "'146.18"
$ws.Range("E25").Value = "'  -2.03%  "
$ws.Range("D26").Value = "'7.23"
$ws.Range("E26").Value = "'  -2.62%  "
$ws.Range("D27").Value = "'16.25"
$ws.Range("E27").Value = "'  -2.08%  "
$ws.Range("E28").Value = "'  -0.14%  "
$ws.Range("E29").Value = "'  -2.03%  "
$ws.Range("E30").Value = "'  -3.03%  "
$ws.Range("E31").Value = "'  -0.58%  "
$ws.Range("D32").Value = "'3.32"
$ws.Range("E32").Value = "'  -2.79%  "
$ws.Range("D33").Value = "'1.455.86"
$ws.Range("E33").Value = "'  -1.95%  "
$ws.Range("E34").Value = "'  -4.28%  "
$ws.Range("E35").Value = "'  -4.32%  "
$ws.Range("D36").Value = "'2.40"
$ws.Range("E36").Value = "'  -0.55%  "
$ws.Range("D37").Value = "'0.921"
$ws.Range("E37").Value = "'  -5.58%  "
$ws.Range("D38").Value = "'0.572"
$ws.Range("E38").Value = "'  -4.50%  "
$ws.Range("E39").Value = "'  -2.99%  "
$ws.Range("E40").Value = "'  +0.27%  "
$ws.Range("E41").Value = "'  -0.03%  "
$ws.Range("D42").Value = "'66.68"
$ws.Range("E42").Value = "'  -4.19%  "
$ws.Range("E43").Value = "'  -3.11%  "
$ws.Range("E44").Value = "'  -2.75%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "'1.803.92"
$ws.Range("E45").Value = "'  -3.77%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'0.790"
$ws.Range("E46").Value = "'  -1.96%  "
$ws.Range("E47").Value = "'  +0.21%  "
$ws.Range("D48").Value = "'88.75"
$ws.Range("E48").Value = "'  -1.74%  "
$ws.Range("E49").Value = "'  -6.05%  "
$ws.Range("E50").Value = "'  -1.45%  "
$ws.Range("D51").Value = "'7.83"
$ws.Range("E51").Value = "'  -3.31%  "
